$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = '@'
$c.Value = '28.089.73'
$c.Style = 'Normal'
$c = $ws.Range('E2')
$c.NumberFormat = '@'
$c.Value = '  -1.62%  '
$c.Style = 'Normal'
$c = $ws.Range('D3')
$c.NumberFormat = '@'
$c.Value = '1.794.21'
$c.Style = 'Normal'
$c = $ws.Range('E3')
$c.NumberFormat = '@'
$c.Value = '  -1.75%  '
$c.Style = 'Normal'
$c = $ws.Range('D4')
$c.NumberFormat = '@'
$c.Value = '1.003'
$c.Style = 'Normal'
$c = $ws.Range('E4')
$c.NumberFormat = '@'
$c.Value = '  +0.25%  '
$c.Style = 'Normal'
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '313.70'
$c.Style = 'Normal'
$c = $ws.Range('E5')
$c.NumberFormat = '@'
$c.Value = '  -0.61%  '
$c.Style = 'Normal'
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '1.002'
$c.Style = 'Normal'
$c = $ws.Range('E6')
$c.NumberFormat = '@'
$c.Value = '  +0.19%  '
$c.Style = 'Normal'
$c = $ws.Range('D7')
$c.NumberFormat = '@'
$c.Value = '0.5207'
$c.Style = 'Normal'
$c = $ws.Range('E7')
$c.NumberFormat = '@'
$c.Value = '  +1.76%  '
$c.Style = 'Normal'
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '0.3821'
$c.Style = 'Normal'
$c = $ws.Range('E8')
$c.NumberFormat = '@'
$c.Value = '  -3.45%  '
$c.Style = 'Normal'
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '0.07878'
$c.Style = 'Normal'
$c = $ws.Range('E9')
$c.NumberFormat = '@'
$c.Value = '  -4.06%  '
$c.Style = 'Normal'
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '41.35'
$c.Style = 'Normal'
$c = $ws.Range('E10')
$c.NumberFormat = '@'
$c.Value = '  -0.92%  '
$c.Style = 'Normal'
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '1.097'
$c.Style = 'Normal'
$c = $ws.Range('E11')
$c.NumberFormat = '@'
$c.Value = '  -1.81%  '
$c.Style = 'Normal'
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '6.276'
$c.Style = 'Normal'
$c = $ws.Range('E12')
$c.NumberFormat = '@'
$c.Value = '  -1.05%  '
$c.Style = 'Normal'
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '1.003'
$c.Style = 'Normal'
$c = $ws.Range('E13')
$c.NumberFormat = '@'
$c.Value = '  +0.24%  '
$c.Style = 'Normal'
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '20.51'
$c.Style = 'Normal'
$c = $ws.Range('E14')
$c.NumberFormat = '@'
$c.Value = '  -3.13%  '
$c.Style = 'Normal'
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '1.795.67'
$c.Style = 'Normal'
$c = $ws.Range('E15')
$c.NumberFormat = '@'
$c.Value = '  -1.66%  '
$c.Style = 'Normal'
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '7.275'
$c.Style = 'Normal'
$c = $ws.Range('E16')
$c.NumberFormat = '@'
$c.Value = '  -3.74%  '
$c.Style = 'Normal'
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '92.39'
$c.Style = 'Normal'
$c = $ws.Range('E17')
$c.NumberFormat = '@'
$c.Value = '  -0.66%  '
$c.Style = 'Normal'
$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '0.00001083'
$c.Style = 'Normal'
$c = $ws.Range('E18')
$c.NumberFormat = '@'
$c.Value = '  -3.94%  '
$c.Style = 'Normal'
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '0.06546'
$c.Style = 'Normal'
$c = $ws.Range('E19')
$c.NumberFormat = '@'
$c.Value = '  -1.73%  '
$c.Style = 'Normal'
$c = $ws.Range('E20')
$c.NumberFormat = '@'
$c.Value = '  +0.27%  '
$c.Style = 'Normal'
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '17.29'
$c.Style = 'Normal'
$c = $ws.Range('E21')
$c.NumberFormat = '@'
$c.Value = '  -3.15%  '
$c.Style = 'Normal'
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '5.950'
$c.Style = 'Normal'
$c = $ws.Range('E22')
$c.NumberFormat = '@'
$c.Value = '  -2.46%  '
$c.Style = 'Normal'
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '28.118.27'
$c.Style = 'Normal'
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '11.12'
$c.Style = 'Normal'
$c = $ws.Range('E24')
$c.NumberFormat = '@'
$c.Value = '  -2.71%  '
$c.Style = 'Normal'
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '2.259'
$c.Style = 'Normal'
$c = $ws.Range('E25')
$c.NumberFormat = '@'
$c.Value = '  -0.29%  '
$c.Style = 'Normal'
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '160.88'
$c.Style = 'Normal'
$c = $ws.Range('E26')
$c.NumberFormat = '@'
$c.Value = '  +2.68%  '
$c.Style = 'Normal'
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '20.46'
$c.Style = 'Normal'
$c = $ws.Range('E27')
$c.NumberFormat = '@'
$c.Value = '  -4.43%  '
$c.Style = 'Normal'
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '1.995.73'
$c.Style = 'Normal'
$c = $ws.Range('E28')
$c.NumberFormat = '@'
$c.Value = '  -1.85%  '
$c.Style = 'Normal'
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '2.333'
$c.Style = 'Normal'
$c = $ws.Range('E29')
$c.NumberFormat = '@'
$c.Value = '  -3.62%  '
$c.Style = 'Normal'
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '122.92'
$c.Style = 'Normal'
$c = $ws.Range('E30')
$c.NumberFormat = '@'
$c.Value = '  -3.34%  '
$c.Style = 'Normal'
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '0.1064'
$c.Style = 'Normal'
$c = $ws.Range('E31')
$c.NumberFormat = '@'
$c.Value = '  -2.29%  '
$c.Style = 'Normal'
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '1.051'
$c.Style = 'Normal'
$c = $ws.Range('E32')
$c.NumberFormat = '@'
$c.Value = '  -5.75%  '
$c.Style = 'Normal'
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '3.672'
$c.Style = 'Normal'
$c = $ws.Range('E33')
$c.NumberFormat = '@'
$c.Value = '  +0.42%  '
$c.Style = 'Normal'
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '5.540'
$c.Style = 'Normal'
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '0.07251'
$c.Style = 'Normal'
$c = $ws.Range('E35')
$c.NumberFormat = '@'
$c.Value = '  +2.60%  '
$c.Style = 'Normal'
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '12.21'
$c.Style = 'Normal'
$c = $ws.Range('E36')
$c.NumberFormat = '@'
$c.Value = '  +8.02%  '
$c.Style = 'Normal'
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '0.02315'
$c.Style = 'Normal'
$c = $ws.Range('E37')
$c.NumberFormat = '@'
$c.Value = '  -1.75%  '
$c.Style = 'Normal'
$c = $ws.Range('B38')
$c.NumberFormat = '@'
$c.Value = 'Algorand'
$c.Style = 'Normal'
$c = $ws.Range('C38')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$c.Style = 'Normal'
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '0.2134'
$c.Style = 'Normal'
$c = $ws.Range('E38')
$c.NumberFormat = '@'
$c.Value = '  -4.40%  '
$c.Style = 'Normal'
$c = $ws.Range('B39')
$c.NumberFormat = '@'
$c.Value = 'FraxShare'
$c.Style = 'Normal'
$c = $ws.Range('C39')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$c.Style = 'Normal'
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '8.719'
$c.Style = 'Normal'
$c = $ws.Range('E39')
$c.NumberFormat = '@'
$c.Value = '  -1.20%  '
$c.Style = 'Normal'
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '5.056'
$c.Style = 'Normal'
$c = $ws.Range('E40')
$c.NumberFormat = '@'
$c.Value = '  -4.41%  '
$c.Style = 'Normal'
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '0.6140'
$c.Style = 'Normal'
$c = $ws.Range('E41')
$c.NumberFormat = '@'
$c.Value = '  -3.12%  '
$c.Style = 'Normal'
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '1.160'
$c.Style = 'Normal'
$c = $ws.Range('E42')
$c.NumberFormat = '@'
$c.Value = '  -2.05%  '
$c.Style = 'Normal'
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '1.376'
$c.Style = 'Normal'
$c = $ws.Range('E43')
$c.NumberFormat = '@'
$c.Value = '  -1.57%  '
$c.Style = 'Normal'
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '13.16'
$c.Style = 'Normal'
$c = $ws.Range('E44')
$c.NumberFormat = '@'
$c.Value = '  -3.36%  '
$c.Style = 'Normal'
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '0.5929'
$c.Style = 'Normal'
$c = $ws.Range('E46')
$c.NumberFormat = '@'
$c.Value = '  -0.42%  '
$c.Style = 'Normal'
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '127.70'
$c.Style = 'Normal'
$c = $ws.Range('E47')
$c.NumberFormat = '@'
$c.Value = '  +1.92%  '
$c.Style = 'Normal'
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '1.230'
$c.Style = 'Normal'
$c = $ws.Range('E48')
$c.NumberFormat = '@'
$c.Value = '  +3.16%  '
$c.Style = 'Normal'
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '1.912'
$c.Style = 'Normal'
$c = $ws.Range('E49')
$c.NumberFormat = '@'
$c.Value = '  -4.30%  '
$c.Style = 'Normal'
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '0.06738'
$c.Style = 'Normal'
$c = $ws.Range('E50')
$c.NumberFormat = '@'
$c.Value = '  -3.05%  '
$c.Style = 'Normal'
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '72.93'
$c.Style = 'Normal'
$c = $ws.Range('E51')
$c.NumberFormat = '@'
$c.Value = '  -1.75%  '
$c.Style = 'Normal'
